# NapoliBallersHistory.xlsx edit:
#  - record the workbook's new save location (absPath)
#  - add a "points for" draft-order / winnings section (cols L & M)
#  - backfill winnings (cols M56:M75) for the 2015/2016 seasons
#  - append the 2017 draft order (rows 76:87)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- workbook-level: absolute path of the (now moved) source file -----
try {
    $wb.Application.ActiveWorkbook.Path = "C:\Users\Tom\Documents\R\NapoliBallers\"
} catch {
    # Path is read-only via COM in most hosts; ignored if unsupported.
}

# --- header row: new draftOrder / winnings columns ---------------------
$ws.Range("L1").Value = "draftOrder"

# --- 2015 season winnings (rows 56-65, ranks 1-10) ----------------------
$ws.Range("M56").Value = 650
$ws.Range("M57").Value = 400
$ws.Range("M58").Value = 200
$ws.Range("M62").Value = -25
$ws.Range("M63").Value = -50
$ws.Range("M64").Value = -75
$ws.Range("M65").Value = -100

# --- 2016 season winnings (rows 66-75, ranks 1-10) ----------------------
$ws.Range("M66").Value = 650
$ws.Range("M67").Value = 400
$ws.Range("M68").Value = 200
$ws.Range("M72").Value = -25
$ws.Range("M73").Value = -50
$ws.Range("M74").Value = -75
$ws.Range("M75").Value = -100

# --- 2017 draft order (new rows 76-87) ----------------------------------
$ws.Range("A76").Value = 2017
$ws.Range("C76").Value = "Terrance Surbella"
$ws.Range("L76").Value = 5

$ws.Range("A77").Value = 2017
$ws.Range("C77").Value = "Alec Emmert"
$ws.Range("L77").Value = 2

$ws.Range("A78").Value = 2017
$ws.Range("C78").Value = "Jake Granatino"
$ws.Range("L78").Value = 12

$ws.Range("A79").Value = 2017
$ws.Range("C79").Value = "Tom Gardner"
$ws.Range("L79").Value = 10

$ws.Range("A80").Value = 2017
$ws.Range("C80").Value = "Tim Lindsay"
$ws.Range("L80").Value = 4

$ws.Range("A81").Value = 2017
$ws.Range("C81").Value = "Hans Biebl"
$ws.Range("L81").Value = 8

$ws.Range("A82").Value = 2017
$ws.Range("C82").Value = "Jesse Burson"
$ws.Range("L82").Value = 9

$ws.Range("A83").Value = 2017
$ws.Range("C83").Value = "william sheridan"
$ws.Range("L83").Value = 6

$ws.Range("A84").Value = 2017
$ws.Range("C84").Value = "Tom Digan"
$ws.Range("L84").Value = 1

$ws.Range("A85").Value = 2017
$ws.Range("C85").Value = "Aaron Peterson"
$ws.Range("L85").Value = 3

$ws.Range("A86").Value = 2017
$ws.Range("C86").Value = "Maksudul Ali"
$ws.Range("L86").Value = 7

$ws.Range("A87").Value = 2017
$ws.Range("C87").Value = "John Ross"
$ws.Range("L87").Value = 11

# --- header row: winnings column (added after the new names so the
#     shared-string table order matches: draftOrder, Maksudul Ali,
#     John Ross, winnings) -------------------------------------------
$ws.Range("M1").Value = "winnings"

# --- column width for the new draftOrder column -------------------------
$ws.Columns.Item(12).ColumnWidth = 10.42578125

# --- view state: scroll + selection left where the author was working --
$win = $excel.ActiveWindow
$win.ScrollRow = 67
$win.ScrollColumn = 1
$ws.Range("M66").Select() | Out-Null
